# Apply edits described by the commit:
#  1. Changed sequence times to 1 second intervals for testing
#     (column A, rows 12-52 on sheet "v1", the first sheet)
#  2. Update selected cell on that sheet from B16 to B11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v1")

# Update the time sequence in column A for rows 12 through 52 so that it
# continues the 1-second-per-row pattern already used in rows 4-11
# (value = (row - 4) seconds expressed as a fraction of a day).
for ($row = 12; $row -le 52; $row++) {
    $seconds = $row - 4
    $ws.Cells.Item($row, 1).Value = $seconds / 86400
}

# Update the active selection on the sheet.
$ws.Activate()
$ws.Range("B11").Select()
